# 1) Split "{issue.bank_commission} руб." into three runs so the
#    templated field becomes "{issue.auto_bank_commission} руб."
#    while keeping identical run formatting (Times New Roman, underline).
$d = $word.ActiveDocument

$hit = $d.Content
$found = $hit.Find.Execute("bank_commission", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'bank_commission' text"
}

$hitStart = $hit.Start

# Insert "auto_" right before "bank_commission"
$insertPoint = $d.Range($hitStart, $hitStart)
$insertPoint.InsertBefore("auto_")

# Force "auto_" to live in its own run (same visible formatting) so the
# final text is split across three <w:r> elements: "{issue.", "auto_",
# "bank_commission} руб." - toggling a property and reverting it keeps
# the run boundary without changing the rendered formatting.
$autoRange = $d.Range($hitStart, $hitStart + 5)
$autoRange.Bold = 1
$autoRange.Bold = 0

# 2) Drop one of the two duplicate trailing empty right-aligned
#    paragraphs, and fold the very last (differently-formatted) empty
#    paragraph into that same right-aligned / zero-spacing style.
$n = $d.Paragraphs.Count
$pDup = $d.Paragraphs.Item($n - 2)
$pDup.Range.Delete()

$n2 = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($n2)
$lastRange = $pLast.Range
$xmlFrag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
           '<w:pPr><w:pStyle w:val="Normal"/><w:spacing w:before="0" w:after="0"/>' + `
           '<w:jc w:val="right"/><w:rPr/></w:pPr><w:r><w:rPr/></w:r></w:p>'
$null = $lastRange.InsertXML($xmlFrag)

$pLastFixed = $d.Paragraphs.Item($d.Paragraphs.Count)
$pLastFixed.Format.SpaceBefore = 0
$pLastFixed.Format.SpaceAfter = 0
